# Project Progress.xlsx - "Proposal Documentation - 2"
#
# The "Module creation" sheet had a spare blank column (old column E, sitting
# between the "ER" and "Select" headers) removed, all of the remaining
# columns were re-sized to a tighter custom layout, and the view was
# scrolled down with B2:J29 left selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the blank spacer column (old column E) - everything to its right
# (Select/Insert/Update/Delete/Print Table/Row Print) shifts one column left.
$ws.Columns.Item(5).Delete()

# Re-apply the custom column widths for the new layout (B..I; J already
# carries over the old column K width unchanged).
$ws.Columns.Item(2).ColumnWidth = 3.5924479166666665   # B
$ws.Columns.Item(3).ColumnWidth = 20.877604166666668   # C
$ws.Columns.Item(4).ColumnWidth = 6.592447916666667    # D
$ws.Columns.Item(5).ColumnWidth = 7.166666666666667    # E
$ws.Columns.Item(6).ColumnWidth = 6.736979166666667    # F
$ws.Columns.Item(7).ColumnWidth = 7.592447916666667    # G
$ws.Columns.Item(8).ColumnWidth = 7.592447916666667    # H
$ws.Columns.Item(9).ColumnWidth = 11.022135416666666   # I

# Scroll the view down (row 16 at the top) and leave the whole table
# (B2:J29) selected, matching the saved view state.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 1
$ws.Range("B2:J29").Select()
